$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{Cell="D2"; Value="26.923.74"},
    @{Cell="E2"; Value="  -2.01%  "},
    @{Cell="D3"; Value="1.834.46"},
    @{Cell="E3"; Value="  -1.63%  "},
    @{Cell="D4"; Value="1.005"},
    @{Cell="E4"; Value="  +0.04%  "},
    @{Cell="D5"; Value="310.52"},
    @{Cell="E5"; Value="  -1.71%  "},
    @{Cell="E6"; Value="  +0.05%  "},
    @{Cell="D7"; Value="0.4620"},
    @{Cell="E7"; Value="  -0.88%  "},
    @{Cell="D8"; Value="0.3660"},
    @{Cell="E8"; Value="  -1.91%  "},
    @{Cell="D9"; Value="0.07173"},
    @{Cell="E9"; Value="  -2.90%  "},
    @{Cell="D10"; Value="0.8804"},
    @{Cell="E10"; Value="  -0.87%  "},
    @{Cell="D11"; Value="0.07834"},
    @{Cell="E11"; Value="  -1.59%  "},
    @{Cell="D12"; Value="19.63"},
    @{Cell="E12"; Value="  -1.94%  "},
    @{Cell="D13"; Value="1.842.15"},
    @{Cell="E13"; Value="  -0.64%  "},
    @{Cell="D14"; Value="5.337"},
    @{Cell="E14"; Value="  -1.71%  "},
    @{Cell="D15"; Value="6.367"},
    @{Cell="E15"; Value="  -3.55%  "},
    @{Cell="D16"; Value="88.70"},
    @{Cell="E16"; Value="  -4.33%  "},
    @{Cell="D17"; Value="1.006"},
    @{Cell="E17"; Value="  -0.04%  "},
    @{Cell="D18"; Value="0.000008759"},
    @{Cell="E18"; Value="  -2.19%  "},
    @{Cell="E19"; Value="  +0.10%  "},
    @{Cell="D20"; Value="26.949.68"},
    @{Cell="E20"; Value="  -2.01%  "},
    @{Cell="D21"; Value="14.51"},
    @{Cell="E21"; Value="  -2.93%  "},
    @{Cell="D22"; Value="5.007"},
    @{Cell="E22"; Value="  -3.35%  "},
    @{Cell="D23"; Value="10.43"},
    @{Cell="E23"; Value="  -1.37%  "},
    @{Cell="E24"; Value="  +4.57%  "},
    @{Cell="D25"; Value="150.77"},
    @{Cell="E25"; Value="  -1.53%  "},
    @{Cell="D26"; Value="18.21"},
    @{Cell="E26"; Value="  -1.70%  "},
    @{Cell="D27"; Value="1.999"},
    @{Cell="E27"; Value="  -4.39%  "},
    @{Cell="D28"; Value="113.69"},
    @{Cell="E28"; Value="  -3.03%  "},
    @{Cell="D29"; Value="4.950"},
    @{Cell="E29"; Value="  -4.13%  "},
    @{Cell="D30"; Value="0.08831"},
    @{Cell="E30"; Value="  -0.91%  "},
    @{Cell="E31"; Value="  +2.98%  "},
    @{Cell="D32"; Value="0.7596"},
    @{Cell="E32"; Value="  +1.10%  "},
    @{Cell="D33"; Value="4.464"},
    @{Cell="E33"; Value="  -0.61%  "},
    @{Cell="D34"; Value="1.137"},
    @{Cell="E34"; Value="  -1.96%  "},
    @{Cell="D35"; Value="2.661"},
    @{Cell="E35"; Value="  +0.27%  "},
    @{Cell="D36"; Value="1.093"},
    @{Cell="E36"; Value="  +1.06%  "},
    @{Cell="D37"; Value="0.01924"},
    @{Cell="E37"; Value="  -2.32%  "},
    @{Cell="D38"; Value="2.930"},
    @{Cell="E38"; Value="  -1.88%  "},
    @{Cell="D39"; Value="0.05124"},
    @{Cell="E39"; Value="  -2.97%  "},
    @{Cell="D40"; Value="6.954"},
    @{Cell="E40"; Value="  -3.25%  "},
    @{Cell="D41"; Value="0.4982"},
    @{Cell="E41"; Value="  -4.76%  "},
    @{Cell="D42"; Value="0.1597"},
    @{Cell="E42"; Value="  -2.93%  "},
    @{Cell="D43"; Value="8.373"},
    @{Cell="E43"; Value="  +0.19%  "},
    @{Cell="B44"; Value="Decentraland"},
    @{Cell="C44"; Value="https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"},
    @{Cell="D44"; Value="0.4673"},
    @{Cell="E44"; Value="  -4.52%  "},
    @{Cell="B45"; Value="PaxDollar"},
    @{Cell="C45"; Value="https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"},
    @{Cell="D45"; Value="1.005"},
    @{Cell="E45"; Value="  +0.06%  "},
    @{Cell="B46"; Value="EnergySwap"},
    @{Cell="C46"; Value="https://coinranking.com/coin/SbWqqTui-+energyswap-ens"},
    @{Cell="D46"; Value="10.19"},
    @{Cell="E46"; Value="  -1.68%  "},
    @{Cell="D47"; Value="102.65"},
    @{Cell="E47"; Value="  -0.93%  "},
    @{Cell="D48"; Value="1.613"},
    @{Cell="E48"; Value="  -3.40%  "},
    @{Cell="D49"; Value="0.06097"},
    @{Cell="E49"; Value="  -2.64%  "},
    @{Cell="D50"; Value="64.66"},
    @{Cell="E50"; Value="  -1.96%  "},
    @{Cell="D51"; Value="36.42"},
    @{Cell="E51"; Value="  -2.38%  "}
)

foreach ($chg in $changes) {
    $rng = $ws.Range($chg.Cell)
    if ($chg.Cell[0] -eq "D" -or $chg.Cell[0] -eq "E") {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $chg.Value
}